$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6681570004204502
$ws.Range("C2").Value = 0.1604053655958566
$ws.Range("D2").Value = 0.04564022844408555
$ws.Range("E2").Value = 0.1136851970799277
$ws.Range("F2").Value = 0.9761275492257937
$ws.Range("K2").Value = 0.3453905516444422
$ws.Range("L2").Value = 0.194925035562008
$ws.Range("N2").Value = 1.777081668538673
$ws.Range("O2").Value = 3.483441505578298
$ws.Range("B3").Value = 0.625680176426755
$ws.Range("C3").Value = 0.1602524579799933
$ws.Range("D3").Value = 0.04377703379804387
$ws.Range("E3").Value = 0.1131101741653424
$ws.Range("F3").Value = 0.9753539534409583
$ws.Range("K3").Value = 0.308724455370168
$ws.Range("L3").Value = 0.1878130467409278
$ws.Range("N3").Value = 1.795260466066001
$ws.Range("O3").Value = 3.494724656376093
$ws.Range("B4").Value = 0.599842645944932
$ws.Range("C4").Value = 0.16016898756369
$ws.Range("D4").Value = 0.04261964192379963
$ws.Range("E4").Value = 0.112813985447616
$ws.Range("F4").Value = 0.9754090563247715
$ws.Range("K4").Value = 0.2862491978763728
$ws.Range("L4").Value = 0.1835499510698924
$ws.Range("N4").Value = 1.806992143303773
$ws.Range("O4").Value = 3.503407944328444
$ws.Range("B5").Value = 0.5893754893466792
$ws.Range("C5").Value = 0.160137611736026
$ws.Range("D5").Value = 0.04214465020783109
$ws.Range("E5").Value = 0.1127076087529844
$ws.Range("F5").Value = 0.9755648813239688
$ws.Range("K5").Value = 0.2771003453930376
$ws.Range("L5").Value = 0.1818388534495909
$ws.Range("N5").Value = 1.811916181421172
$ws.Range("O5").Value = 3.507388020435627
$ws.Range("B6").Value = 0.5876411801841073
$ws.Range("C6").Value = 0.1601325617832927
$ws.Range("D6").Value = 0.04206557674890377
$ws.Range("E6").Value = 0.112690810630788
$ws.Range("F6").Value = 0.97559881335075
$ws.Range("K6").Value = 0.2755818031805006
$ws.Range("L6").Value = 0.1815563084687142
$ws.Range("N6").Value = 1.81274246820559
$ws.Range("O6").Value = 3.50807558325306
$ws.Range("B7").Value = 0.5997012309476304
$ws.Range("C7").Value = 0.1601685537065514
$ws.Range("D7").Value = 0.04261324953650103
$ws.Range("E7").Value = 0.1128124927954381
$ws.Range("F7").Value = 0.9754106177275759
$ws.Range("K7").Value = 0.2861257721451409
$ws.Range("L7").Value = 0.1835267686367672
$ws.Range("N7").Value = 1.807057970453685
$ws.Range("O7").Value = 3.503459832935448
$ws.Range("B8").Value = 0.653460840824124
$ws.Range("C8").Value = 0.1603504931666144
$ws.Range("D8").Value = 0.04500058810160823
$ws.Range("E8").Value = 0.1134751359154507
$ws.Range("F8").Value = 0.9757508337641383
$ws.Range("K8").Value = 0.3327404938338248
$ws.Range("L8").Value = 0.1924513417285567
$ws.Range("N8").Value = 1.783231451866865
$ws.Range("O8").Value = 3.48696769487924
$ws.Range("B9").Value = 0.7607922397304776
$ws.Range("C9").Value = 0.160789084671066
$ws.Range("D9").Value = 0.04957527409242601
$ws.Range("E9").Value = 0.1152251885736355
$ws.Range("F9").Value = 0.9806221767804644
$ws.Range("K9").Value = 0.4244360343207063
$ws.Range("L9").Value = 0.2107730314315006
$ws.Range("N9").Value = 1.741028368442693
$ws.Range("O9").Value = 3.468550780601646
$ws.Range("B10").Value = 0.8407906717435765
$ws.Range("C10").Value = 0.1611601358179726
$ws.Range("D10").Value = 0.05287049434147661
$ws.Range("E10").Value = 0.1167850329405873
$ws.Range("F10").Value = 0.9867636300128311
$ws.Range("K10").Value = 0.4919627630806076
$ws.Range("L10").Value = 0.2247330074573171
$ws.Range("N10").Value = 1.712775750114063
$ws.Range("O10").Value = 3.463507185780543
$ws.Range("B11").Value = 0.8774279078903646
$ws.Range("C11").Value = 0.1613393189855969
$ws.Range("D11").Value = 0.05435516573734134
$ws.Range("E11").Value = 0.117554046804738
$ws.Range("F11").Value = 0.9901139800157068
$ws.Range("K11").Value = 0.5227138875153798
$ws.Range("L11").Value = 0.2311919513844032
$ws.Range("N11").Value = 1.700520497415737
$ws.Range("O11").Value = 3.463055735570066
$ws.Range("B12").Value = 0.8913362473992947
$ws.Range("C12").Value = 0.1614086459381951
$ws.Range("D12").Value = 0.05491529280831031
$ws.Range("E12").Value = 0.1178537829382371
$ws.Range("F12").Value = 0.9914626671603912
$ws.Range("K12").Value = 0.5343628709390771
$ws.Range("L12").Value = 0.2336533325760257
$ws.Range("N12").Value = 1.695965610553349
$ws.Range("O12").Value = 3.463149727143389
$ws.Range("B13").Value = 0.8883393092922915
$ws.Range("C13").Value = 0.1613936498909254
$ws.Range("D13").Value = 0.05479475256210975
$ws.Range("E13").Value = 0.1177888504513263
$ws.Range("F13").Value = 0.9911686468553995
$ws.Range("K13").Value = 0.5318538755439874
$ws.Range("L13").Value = 0.2331225411529374
$ws.Range("N13").Value = 1.696942764978462
$ws.Range("O13").Value = 3.463117701820977
$ws.Range("B14").Value = 0.8785714656414143
$ws.Range("C14").Value = 0.1613449931249846
$ws.Range("D14").Value = 0.0544012896655488
$ws.Range("E14").Value = 0.1175785354873931
$ws.Range("F14").Value = 0.9902233346552691
$ws.Range("K14").Value = 0.523672175242865
$ws.Range("L14").Value = 0.2313941401270654
$ws.Range("N14").Value = 1.700144041971894
$ws.Range("O14").Value = 3.463058158756809
$ws.Range("B15").Value = 0.8725928663614866
$ws.Range("C15").Value = 0.1613153808436252
$ws.Range("D15").Value = 0.05416000995204939
$ws.Range("E15").Value = 0.1174508214333336
$ws.Range("F15").Value = 0.9896547181361655
$ws.Range("K15").Value = 0.5186611788586504
$ws.Range("L15").Value = 0.2303374632766264
$ws.Range("N15").Value = 1.702116106667567
$ws.Range("O15").Value = 3.463056188856058
$ws.Range("B16").Value = 0.8384012237600871
$ws.Range("C16").Value = 0.1611486331473699
$ws.Range("D16").Value = 0.05277317698701012
$ws.Range("E16").Value = 0.1167359703864257
$ws.Range("F16").Value = 0.9865558719802721
$ws.Range("K16").Value = 0.4899537212522205
$ws.Range("L16").Value = 0.2243130763585981
$ws.Range("N16").Value = 1.713588678273044
$ws.Range("O16").Value = 3.46357376783655
$ws.Range("B17").Value = 0.8174881541653747
$ws.Range("C17").Value = 0.1610489857895487
$ws.Range("D17").Value = 0.05191871009629523
$ws.Range("E17").Value = 0.1163126410416986
$ws.Range("F17").Value = 0.9847973448938063
$ws.Range("K17").Value = 0.4723506711353593
$ws.Range("L17").Value = 0.2206450378355385
$ws.Range("N17").Value = 1.720779680501872
$ws.Range("O17").Value = 3.464363254555138
$ws.Range("B18").Value = 0.8054826634906931
$ws.Range("C18").Value = 0.1609926504160484
$ws.Range("D18").Value = 0.0514258949258064
$ws.Range("E18").Value = 0.1160747484478115
$ws.Range("F18").Value = 0.9838382769862264
$ws.Range("K18").Value = 0.4622289883931785
$ws.Range("L18").Value = 0.2185454946470458
$ws.Range("N18").Value = 1.724971943806412
$ws.Range("O18").Value = 3.464990813210079
$ws.Range("B19").Value = 0.8014218085243101
$ws.Range("C19").Value = 0.1609737450063022
$ws.Range("D19").Value = 0.05125880530331983
$ws.Range("E19").Value = 0.1159951636724834
$ws.Range("F19").Value = 0.9835225527429756
$ws.Range("K19").Value = 0.4588025135104772
$ws.Range("L19").Value = 0.2178363823608578
$ws.Range("N19").Value = 1.726401020780456
$ws.Range("O19").Value = 3.46523308875345
$ws.Range("B20").Value = 0.8197119948320335
$ws.Range("C20").Value = 0.1610594922270181
$ws.Range("D20").Value = 0.05200980928198362
$ws.Range("E20").Value = 0.1163571262241625
$ws.Range("F20").Value = 0.9849791212725449
$ws.Range("K20").Value = 0.4742242261522165
$ws.Range("L20").Value = 0.2210344501056625
$ws.Range("N20").Value = 1.720008370237489
$ws.Range("O20").Value = 3.464261259627506
$ws.Range("B21").Value = 0.8814395842542808
$ws.Range("C21").Value = 0.1613592449455936
$ws.Range("D21").Value = 0.05451691605161813
$ws.Range("E21").Value = 0.1176400788153131
$ws.Range("F21").Value = 0.9904988254851901
$ws.Range("K21").Value = 0.5260752291706581
$ws.Range("L21").Value = 0.2319013928576084
$ws.Range("N21").Value = 1.69920141719581
$ws.Range("O21").Value = 3.463068457935037
$ws.Range("B22").Value = 0.921983561466277
$ws.Range("C22").Value = 0.1615637341537877
$ws.Range("D22").Value = 0.05614328684723091
$ws.Range("E22").Value = 0.1185282570388715
$ws.Range("F22").Value = 0.9945724446693731
$ws.Range("K22").Value = 0.5599870924885693
$ws.Range("L22").Value = 0.2390939889457115
$ws.Range("N22").Value = 1.686103658704706
$ws.Range("O22").Value = 3.463833171941559
$ws.Range("B23").Value = 0.9003262772671405
$ws.Range("C23").Value = 0.1614538155716332
$ws.Range("D23").Value = 0.05527638335661322
$ws.Range("E23").Value = 0.1180496788865746
$ws.Range("F23").Value = 0.9923556399580065
$ws.Range("K23").Value = 0.5418856497705349
$ws.Range("L23").Value = 0.2352469181527255
$ws.Range("N23").Value = 1.693048335354891
$ws.Range("O23").Value = 3.463283748597632
$ws.Range("B24").Value = 0.8187065412442678
$ws.Range("C24").Value = 0.1610547392956718
$ws.Range("D24").Value = 0.05196862823077453
$ws.Range("E24").Value = 0.1163369973866359
$ws.Range("F24").Value = 0.9848967784005112
$ws.Range("K24").Value = 0.4733771964216089
$ws.Range("L24").Value = 0.2208583679661729
$ws.Range("N24").Value = 1.720356898755047
$ws.Range("O24").Value = 3.464306830548281
$ws.Range("B25").Value = 0.7315541259502254
$ws.Range("C25").Value = 0.160661785395888
$ws.Range("D25").Value = 0.04834920894071359
$ws.Range("E25").Value = 0.1147035761417037
$ws.Range("F25").Value = 0.9788543130697533
$ws.Range("K25").Value = 0.3996011089199101
$ws.Range("L25").Value = 0.2057288020059787
$ws.Range("N25").Value = 1.751961663483082
$ws.Range("O25").Value = 3.472042496418453
